$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("predictions")

# Rename the ambiguous prediction-table headers for clarity.
# Writing directly into the header row keeps the underlying ListObject
# (Table2) column names in sync.
$ws.Range("D1").Value = "team1wins"
$ws.Range("E1").Value = "team2wins"
$ws.Range("F1").Value = "draw"
$ws.Range("I1").Value = "team1won"
$ws.Range("J1").Value = "team2wins2"
$ws.Range("K1").Value = "noteamwins"
$ws.Range("O1").Value = "loss"

# The structured-reference formulas in the calculated columns refer to the
# old header names, so re-point them at the new ones (same calculation,
# new column labels). Set row-by-row so each cell keeps its own formula
# (matching how the workbook already stores them) instead of becoming one
# shared-formula block.
for ($r = 2; $r -le 13; $r++) {
    $ws.Range("L$r").Formula = "=-LOG(Table2[[#This Row],[team1wins]], EXP(1))*Table2[[#This Row],[team1won]]"
    $ws.Range("M$r").Formula = "=-LOG(Table2[[#This Row],[team2wins]], EXP(1))*Table2[[#This Row],[team2wins2]]"
    $ws.Range("N$r").Formula = "=-LOG(Table2[[#This Row],[draw]], EXP(1))*Table2[[#This Row],[noteamwins]]"
}

$ws.Range("I20").Select()

$wb.Save()
